$d = $word.ActiveDocument

# Find the paragraph that contains the "LOQ4083 ... (Requisito)" line.
# Right after it the document currently has four extra paragraphs that
# need to go away:
#   1. an empty paragraph
#   2. "Ver no Jupiter Salvar em pdf Salvar em docx"
#   3. an empty paragraph
#   4. an empty paragraph with a page break before it (and jc=left)
# Everything from paragraph 5 onward (an empty paragraph followed by a
# page-break paragraph, then the section properties) must be preserved.
$searchRange = $d.Content
$found = $searchRange.Find.Execute("LOQ4083", $true, $false, $false, $false, $false, `
                                    $true, 1, $false, "", 0)

if ($found) {
    $anchorIndex = $searchRange.Paragraphs.First.Index

    $startPara = $d.Paragraphs.Item($anchorIndex + 1)
    $endPara = $d.Paragraphs.Item($anchorIndex + 5)

    $deleteRange = $d.Range($startPara.Range.Start, $endPara.Range.Start)
    $deleteRange.Delete()
}
